$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "You have a tensor y of shape 1234x25 and you want to refactor it so that each row will contain 50 values. Which operation should you use to accomplish this?",
        "ques_type": 2,
        "options": [
            "y.shape(617, 50)",
            "y.reshape(1, -2)",
            "y.view(-1, 50)",
            "y.transform(1, 2)"
        ],
        "score": "y.view(-1, 50)"
    },
    {
        "title": "You want to apply data augmentation on a dataset with images. The images have already been converted to tensors. Now you just want to crop the tensors representing the images from the corners and from the center, thus augmenting several new images from each tensor. The augmented images will not take any other modifications. Which built-in transform class in PyTorch should you use?",
        "ques_type": 2,
        "options": [
            "RandomCrop",
            "CornerCrop",
            "FiveCrop",
            "CenterCrop"
        ],
        "score": "FiveCrop"
    },
    {
        "title": "You have created a MaxPool layer by executing the following code: nn.MaxPool1d(1, stride=2) Now you want to test the execution of the layer by providing a tensor of shape (4,2) as an input to it. What will be the shape of the output tensor?",
        "ques_type": 2,
        "options": [
            "(4, 0)",
            "(2, 1)",
            "(4, 1)",
            "(2, 0)"
        ],
        "score": "(4, 1)"
    },
    {
        "title": "You have a tensor of shape 10x20x30, and you want to flatten it into a two-dimensional tensor of shape 200x30 to use it as a layer in a neural network. You want to use the built-in Flatten class available in PyTorch. How should you initialize the class?",
        "ques_type": 2,
        "options": [
            "Flatten()",
            "Flatten(0)",
            "Flatten(0, 1)",
            "Flatten(0, -1)"
        ],
        "score": "Flatten(0, 1)"
    }
]
'@

$ws.Range("A2").ClearContents()
$ws.Range("A1").Value = $text
$ws.Range("A1").Style = "Normal"
$ws.Rows.Item(1).EntireRow.AutoFit()
